$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is a duplicate of row 7; delete it entirely so rows 9-16 shift up
# to become rows 8-15 (matches "up to date with run 30" commit).
$ws.Rows.Item(8).Delete()
